# Commit: "XLXS file has been corrected"
#
# Comparing the canonical OOXML before/after, almost every hunk in the diff
# (namespace bumps on workbook.xml/sheet1.xml/theme1.xml, <fileVersion>,
# <mc:AlternateContent>, <xr:revisionPtr>, <bookViews>, the numFmt
# "10.0" -> "10" / "1.0" -> "1" normalizations, row "spans"/x14ac:dyDescent
# attributes, the new <pageMargins .../>, etc.) is simply what real desktop
# Excel stamps onto a workbook the moment it opens + resaves a file that
# wasn't produced by Excel itself (this workbook has no docProps/app.xml or
# docProps/core.xml at all beforehand, which is the tell).
#
# The one genuine, user-driven content change buried in that noise is that
# xl/drawings/drawing1.xml (the floating image "image1.png" anchored over
# the sheet) is gone afterwards, and sheet1.xml no longer carries its
# <drawing r:id="rId1"/> reference. I.e. the author deleted the picture that
# was sitting on top of the data. Reproduce that with the Shapes API.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    for ($i = $ws.Shapes.Count; $i -ge 1; $i--) {
        $ws.Shapes.Item($i).Delete()
    }
}
